$wb = $excel.ActiveWorkbook
$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll = $wb.Worksheets.Item("全部类型")

$wsExhibit.Range("F2").Value = 7004
$wsExhibit.Range("F4").Value = 460
$wsExhibit.Range("F6").Value = 552
$wsExhibit.Range("F7").Value = 142
$wsExhibit.Range("F8").Value = 118
$wsExhibit.Range("F11").Value = 51
$wsExhibit.Range("F12").Value = 199
$wsExhibit.Range("F13").Value = 443
$wsExhibit.Range("F14").Value = 27
$wsExhibit.Range("F15").Value = 1822
$wsExhibit.Range("F17").Value = 3612
$wsExhibit.Range("F19").Value = 246
$wsExhibit.Range("F20").Value = 83
$wsExhibit.Range("F21").Value = 22
$wsExhibit.Range("F22").Value = 26
$wsExhibit.Range("F23").Value = 2240
$wsExhibit.Range("F24").Value = 15
$wsExhibit.Range("F25").Value = 247
$wsExhibit.Range("F27").Value = 34
$wsExhibit.Range("F31").Value = 158
$wsExhibit.Range("F32").Value = 231
$wsExhibit.Range("G32").Value = 69
$wsExhibit.Range("F33").Value = 86
$wsAll.Range("F2").Value = 7004
$wsAll.Range("F4").Value = 460
$wsAll.Range("F7").Value = 552
$wsAll.Range("F8").Value = 142
$wsAll.Range("F9").Value = 118
$wsAll.Range("F12").Value = 51
$wsAll.Range("F13").Value = 199
$wsAll.Range("F14").Value = 443
$wsAll.Range("F15").Value = 27
$wsAll.Range("F16").Value = 1822
$wsAll.Range("F18").Value = 3612
$wsAll.Range("F20").Value = 246
$wsAll.Range("F21").Value = 83
$wsAll.Range("F22").Value = 22
$wsAll.Range("F23").Value = 26
$wsAll.Range("F24").Value = 2240
$wsAll.Range("F25").Value = 15
$wsAll.Range("F26").Value = 247
$wsAll.Range("F28").Value = 34
$wsAll.Range("F32").Value = 158
$wsAll.Range("F33").Value = 231
$wsAll.Range("G33").Value = 69
$wsAll.Range("F34").Value = 86
